# Regenerate save_data column G (header "K", formerly based on "Strike#")
# with recalculated values after switching the strike-number convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (rows 2-30), per the recalculated "K" metric.
$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 1
    24 = 0
    25 = 2
    26 = 4
    27 = 3
    28 = 1
    29 = 1
    30 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
